# Commit: "add words in March 9th second"
#
# The document has three paragraphs:
#   1. "2023年3月8号"
#   2. "下雨，今天女神节，又是happy的一天，祝福所有女生节日快乐，健康。"
#   3. "2023年3月9号"   (this paragraph owns a <w:pPr> with rFonts hint=eastAsia)
#
# The edit inserts a brand-new paragraph *before* paragraph 3, re-using the
# exact same text/run-split as paragraph 3 ("2023年3月9号" split into three
# runs: "2", "023", "年3月9号" -- the first and third carrying
# <w:rFonts w:hint="eastAsia"/>), but with NO paragraph properties of its
# own.
#
# The original paragraph 3 (which keeps its <w:pPr>) then has its text
# replaced with "晴，天气很好".

$d = $word.ActiveDocument

# Locate the paragraph that currently reads "2023年3月9号".
$found = $d.Content
$found.Find.Execute("2023年3月9号", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
$datePara = $found.Paragraphs(1)

# Work out that paragraph's 1-based index so we can address paragraphs by
# position afterwards (Find-based re-lookup is ambiguous once the text has
# been duplicated into a second paragraph).
$dateIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Start -eq $datePara.Range.Start) {
        $dateIndex = $i
    }
}

# Insert a brand new (empty) paragraph immediately before it; this becomes
# paragraph number $dateIndex, and the original date paragraph shifts to
# $dateIndex + 1.
$datePara.Range.InsertParagraphBefore() | Out-Null

$newPara = $d.Paragraphs($dateIndex)
$oldPara = $d.Paragraphs($dateIndex + 1)

# Fill the new paragraph with the three runs, reproducing the exact
# formatting (eastAsia hint on the CJK-containing runs) via OOXML, and no
# paragraph properties -- matching the target markup precisely.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t>023</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>年3月9号</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($xml) | Out-Null

# The original paragraph (the one that still owns the <w:pPr>) keeps its
# position right after the freshly-created one; change its text to
# "晴，天气很好" while preserving its paragraph/run formatting.
$textRange = $oldPara.Range.Duplicate
$textRange.MoveEnd(1, -1) | Out-Null
$textRange.Text = "晴，天气很好"
